# Auto-generated edit script applying numeric updates described in the commit diff.
# Updates currentAveragePrice / currentAveragePriceHQ / LevePriceHQ / LeveProfitHQ (and related)
# columns (H, I, J, K, L, M, N) for specific Leve rows across all 8 sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 108
$ws.Range("H108").Value = 32148.834
$ws.Range("J108").Value = 32148.834
$ws.Range("L108").Value = 32148.834
$ws.Range("N108").Value = -39828.834
# Row 126
$ws.Range("H126").Value = 42152.5
$ws.Range("J126").Value = 42152.5
$ws.Range("L126").Value = 42152.5
$ws.Range("N126").Value = -52032.5
# Row 128
$ws.Range("H128").Value = 36592
$ws.Range("J128").Value = 36592
$ws.Range("L128").Value = 36592
$ws.Range("N128").Value = -46552
# Row 130
$ws.Range("H130").Value = 42265.453
$ws.Range("J130").Value = 42265.453
$ws.Range("L130").Value = 42265.453
$ws.Range("N130").Value = -52305.453

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 28449
$ws.Range("J44").Value = 28449
$ws.Range("L44").Value = 28449
$ws.Range("N44").Value = -29425
# Row 80
$ws.Range("H80").Value = 38005
$ws.Range("J80").Value = 38005
$ws.Range("L80").Value = 38005
$ws.Range("N80").Value = -40001
# Row 83
$ws.Range("H83").Value = 38005
$ws.Range("J83").Value = 38005
$ws.Range("L83").Value = 114015
$ws.Range("N83").Value = -123999
# Row 123
$ws.Range("H123").Value = 1278250
$ws.Range("J123").Value = 1278250
$ws.Range("L123").Value = 1278250
$ws.Range("N123").Value = -1288050
# Row 125
$ws.Range("H125").Value = 180038130
$ws.Range("J125").Value = 180038130
$ws.Range("L125").Value = 180038130
$ws.Range("N125").Value = -180047970
# Row 127
$ws.Range("H127").Value = 31557.5
$ws.Range("J127").Value = 31557.5
$ws.Range("L127").Value = 31557.5
$ws.Range("N127").Value = -41477.5
# Row 129
$ws.Range("H129").Value = 46699.75
$ws.Range("J129").Value = 46699.75
$ws.Range("L129").Value = 46699.75
$ws.Range("N129").Value = -56699.75
# Row 134
$ws.Range("H134").Value = 39986
$ws.Range("J134").Value = 39986
$ws.Range("L134").Value = 39986
$ws.Range("N134").Value = -50126
# Row 135
$ws.Range("H135").Value = 35826.75
$ws.Range("J135").Value = 35826.75
$ws.Range("L135").Value = 35826.75
$ws.Range("N135").Value = -45966.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 30549.6
$ws.Range("J35").Value = 30549.6
$ws.Range("L35").Value = 30549.6
$ws.Range("N35").Value = -31169.6
# Row 122
$ws.Range("H122").Value = 48630
$ws.Range("J122").Value = 48630
$ws.Range("L122").Value = 48630
$ws.Range("N122").Value = -58430
# Row 124
$ws.Range("H124").Value = 42500
$ws.Range("J124").Value = 42500
$ws.Range("L124").Value = 42500
$ws.Range("N124").Value = -52320
# Row 125
$ws.Range("H125").Value = 52580
$ws.Range("J125").Value = 52580
$ws.Range("L125").Value = 52580
$ws.Range("N125").Value = -62420
# Row 126
$ws.Range("H126").Value = 33401.25
$ws.Range("J126").Value = 33401.25
$ws.Range("L126").Value = 33401.25
$ws.Range("N126").Value = -43281.25
# Row 132
$ws.Range("H132").Value = 33741.25
$ws.Range("J132").Value = 33741.25
$ws.Range("L132").Value = 33741.25
$ws.Range("N132").Value = -43861.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 19885
$ws.Range("J41").Value = 23846.666
$ws.Range("L41").Value = 23846.666
$ws.Range("N41").Value = -24702.666
# Row 50
$ws.Range("H50").Value = 9098
$ws.Range("J50").Value = 9098
$ws.Range("L50").Value = 9098
$ws.Range("N50").Value = -10348
# Row 51
$ws.Range("H51").Value = 9103.6
$ws.Range("J51").Value = 9103.6
$ws.Range("L51").Value = 9103.6
$ws.Range("N51").Value = -10575.6
# Row 60
$ws.Range("H60").Value = 26118.385
$ws.Range("J60").Value = 26118.385
$ws.Range("L60").Value = 26118.385
$ws.Range("N60").Value = -27140.385
# Row 61
$ws.Range("H61").Value = 9103.6
$ws.Range("J61").Value = 9103.6
$ws.Range("L61").Value = 9103.6
$ws.Range("N61").Value = -9799.6
# Row 68
$ws.Range("H68").Value = 17470.8
$ws.Range("J68").Value = 17470.8
$ws.Range("L68").Value = 17470.8
$ws.Range("N68").Value = -18968.8
# Row 71
$ws.Range("H71").Value = 17470.8
$ws.Range("J71").Value = 17470.8
$ws.Range("L71").Value = 52412.39999999999
$ws.Range("N71").Value = -59900.39999999999
# Row 97
$ws.Range("H97").Value = 9890
$ws.Range("J97").Value = 9890
$ws.Range("L97").Value = 9890
$ws.Range("N97").Value = -11872
# Row 109
$ws.Range("H109").Value = 12000
$ws.Range("J109").Value = 12000
$ws.Range("L109").Value = 12000
$ws.Range("N109").Value = -14080
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()  # was -61873.332
# Row 130
$ws.Range("H130").Value = 54827.5
$ws.Range("J130").Value = 54827.5
$ws.Range("L130").Value = 54827.5
$ws.Range("N130").Value = -64867.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 7312.68
$ws.Range("I87").Value = 2888.1333
$ws.Range("J87").Value = 13949.5
$ws.Range("K87").Value = 8664.3999
$ws.Range("L87").Value = 41848.5
$ws.Range("M87").Value = -7416.3999
$ws.Range("N87").Value = -44344.5
# Row 90
$ws.Range("H90").Value = 7312.68
$ws.Range("I90").Value = 2888.1333
$ws.Range("J90").Value = 13949.5
$ws.Range("K90").Value = 25993.1997
$ws.Range("L90").Value = 125545.5
$ws.Range("M90").Value = -19753.1997
$ws.Range("N90").Value = -138025.5
# Row 92
$ws.Range("H92").Value = 1084
$ws.Range("I92").Value = 912
$ws.Range("J92").Value = 1268.2858
$ws.Range("K92").Value = 2736
$ws.Range("L92").Value = 3804.8574
$ws.Range("M92").Value = -1488
$ws.Range("N92").Value = -6300.857400000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 3351.8
$ws.Range("J43").Value = 9072.666999999999
$ws.Range("L43").Value = 9072.666999999999
$ws.Range("N43").Value = -9374.666999999999
# Row 62
$ws.Range("H62").Value = 34992.5
$ws.Range("J62").Value = 34992.5
$ws.Range("L62").Value = 34992.5
$ws.Range("N62").Value = -36364.5
# Row 65
$ws.Range("H65").Value = 34992.5
$ws.Range("J65").Value = 34992.5
$ws.Range("L65").Value = 104977.5
$ws.Range("N65").Value = -111841.5
# Row 82
$ws.Range("H82").Value = 43994
$ws.Range("J82").Value = 43994
$ws.Range("L82").Value = 43994
$ws.Range("N82").Value = -44760
# Row 85
$ws.Range("H85").Value = 43994
$ws.Range("J85").Value = 43994
$ws.Range("L85").Value = 43994
$ws.Range("N85").Value = -46646
# Row 93
$ws.Range("H93").Value = 9295.916999999999
$ws.Range("J93").Value = 9295.916999999999
$ws.Range("L93").Value = 9295.916999999999
$ws.Range("N93").Value = -13039.917
# Row 122
$ws.Range("H122").Value = 2889.0557
$ws.Range("I122").Value = 2853.7334
$ws.Range("J122").Value = 3065.6667
$ws.Range("K122").Value = 8561.200199999999
$ws.Range("L122").Value = 9197.000100000001
$ws.Range("M122").Value = -6111.200199999999
$ws.Range("N122").Value = -14097.0001
# Row 123
$ws.Range("H123").Value = 31420.666
$ws.Range("J123").Value = 31420.666
$ws.Range("L123").Value = 31420.666
$ws.Range("N123").Value = -36320.666
# Row 124
$ws.Range("H124").Value = 32396
$ws.Range("J124").Value = 32396
$ws.Range("L124").Value = 32396
$ws.Range("N124").Value = -42216
# Row 127
$ws.Range("H127").Value = 56326
$ws.Range("J127").Value = 56326
$ws.Range("L127").Value = 56326
$ws.Range("N127").Value = -66246
# Row 128
$ws.Range("H128").Value = 48525
$ws.Range("J128").Value = 48525
$ws.Range("L128").Value = 48525
$ws.Range("N128").Value = -58485
# Row 135
$ws.Range("H135").Value = 59769.23
$ws.Range("J135").Value = 59769.23
$ws.Range("L135").Value = 59769.23
$ws.Range("N135").Value = -69909.23000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 108
$ws.Range("H108").Value = 24284.25
$ws.Range("J108").Value = 24284.25
$ws.Range("L108").Value = 24284.25
$ws.Range("N108").Value = -31964.25
# Row 109
$ws.Range("H109").Value = 18086.334
$ws.Range("I109").Value = 5259
$ws.Range("K109").Value = 5259
$ws.Range("M109").Value = -3872
# Row 123
$ws.Range("H123").Value = 46563.168
$ws.Range("J123").Value = 46563.168
$ws.Range("L123").Value = 46563.168
$ws.Range("N123").Value = -56363.168
# Row 129
$ws.Range("H129").Value = 42429
$ws.Range("J129").Value = 42429
$ws.Range("L129").Value = 42429
$ws.Range("N129").Value = -52429
# Row 133
$ws.Range("H133").Value = 82800
$ws.Range("J133").Value = 82800
$ws.Range("L133").Value = 82800
$ws.Range("N133").Value = -87860
# Row 134
$ws.Range("H134").Value = 47469.832
$ws.Range("J134").Value = 54085.8
$ws.Range("L134").Value = 54085.8
$ws.Range("N134").Value = -64225.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 108
$ws.Range("H108").Value = 29880
$ws.Range("J108").Value = 29880
$ws.Range("L108").Value = 29880
$ws.Range("N108").Value = -37560
# Row 109
$ws.Range("H109").Value = 19573.5
$ws.Range("J109").Value = 19573.5
$ws.Range("L109").Value = 19573.5
$ws.Range("N109").Value = -22347.5
# Row 127
$ws.Range("H127").Value = 26607.8
$ws.Range("J127").Value = 26607.8
$ws.Range("L127").Value = 26607.8
$ws.Range("N127").Value = -36527.8
# Row 128
$ws.Range("H128").Value = 180701.42
$ws.Range("J128").Value = 180701.42
$ws.Range("L128").Value = 180701.42
$ws.Range("N128").Value = -190661.42
# Row 130
$ws.Range("H130").Value = 24714.5
$ws.Range("J130").Value = 24714.5
$ws.Range("L130").Value = 24714.5
$ws.Range("N130").Value = -34754.5
